$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A data correction (rows 40-47) ---
$ws.Range("A40").Value = 229.21
$ws.Range("A41").Value = 229.21

$ws.Range("A42:A47").Value = 229.226
$ws.Range("A42:A47").HorizontalAlignment = -4108

# --- View: zoom in and move the selection down to the bottom of the table ---
$excel.ActiveWindow.Zoom = 177
$null = $ws.Range("C49").Select()
